$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the sheet view: scroll position + current selection ---
$ws.Activate()
$window = $excel.ActiveWindow
$window.ScrollRow = 8
$window.ScrollColumn = 1
$ws.Range("D12:D17").Select() | Out-Null
$excel.ActiveCell = $ws.Range("D12")

# --- "Banco de Dados / Modelo Lógico" and "Banco de Dados / Relacionamentos"
#     move from "Em andamento" to "Concluído" (reuse the same look already
#     used by the other "Concluído" rows, e.g. D8) ---
$ws.Range("D8").Copy() | Out-Null
$ws.Range("D9:D10").PasteSpecial(-4122) | Out-Null
$ws.Range("D9").Value = "Conclu" + [char]0x00ED + "do"
$ws.Range("D10").Value = "Conclu" + [char]0x00ED + "do"

# --- Restyle the still-"Pendente" backlog rows (Comandos SQL + every
#     Algoritmos sub-item) to match the style already used for D20 ---
$ws.Range("D20").Copy() | Out-Null
$ws.Range("D11:D17").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
